$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

# Populate the "Result" sheet with test run outcomes (rows 2-9), mirroring
# the "Tests" sheet's WorkflowFile/ExpectedResult pairs plus the actual
# Status (PASS/FAIL) and any Comments.
$data = @(
    @("Framework\InitAllSettings.xaml",      "Success",          "PASS", ""),
    @("Framework\InitAllApplications.xaml",  "Success",          "PASS", ""),
    @("Framework\CloseAllApplications.xaml", "Success",          "PASS", ""),
    @("Framework\CloseAllApplications.xaml", "SystemException",  "FAIL", "No exception thrown."),
    @("Framework\InitAllSettings.xaml",      "Success",          "PASS", ""),
    @("Framework\InitAllSettings.xaml",      "Success",          "PASS", ""),
    @("Framework\InitAllApplications.xaml",  "Success",          "PASS", ""),
    @("Framework\CloseAllApplications.xaml", "Success",          "PASS", "")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    if ($entry[3] -ne "") {
        $ws.Cells.Item($row, 4).Value = $entry[3]
    }
    $row++
}

# The Result sheet becomes the active/selected sheet after the run.
$ws.Activate()
